$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.41%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.51"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.51%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.996"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.02%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07808"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.82%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.211"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.75%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.004"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.01%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.991"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.24%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9126"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.96%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09293"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.09%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1858"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.64%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08404"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.06%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03517"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.98%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09936"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.35%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001466"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.19%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005649"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.95%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.476"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.46%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.097"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.29%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.85%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.558"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.90%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2226"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.30%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04633"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.20%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.69%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004441"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.52%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001296"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.40%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004735"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "39.56%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01757"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.94%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04683"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.36%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007831"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.01%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1387"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.71%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007649"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.62%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002292"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.55%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01016"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.66%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006062"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.91%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.43%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.665"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "182.90%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.80%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002093"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.43%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001994"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.43%"
